$wb = $excel.ActiveWorkbook

# --- Sheet 1: ATS Accuracy ---
$ws1 = $wb.Worksheets.Item("ATS Accuracy")

$ws1.Range("B2").Value = 5
$ws1.Range("C2").Value = 86
$ws1.Range("D2").Value = 91
$ws1.Range("E2").Value = 94.5

$ws1.Range("B3").Value = 4
$ws1.Range("C3").Value = 66
$ws1.Range("D3").Value = 70
$ws1.Range("E3").Value = 94.3

$ws1.Range("C4").Value = 9
$ws1.Range("D4").Value = 13
$ws1.Range("E4").Value = 69.2

$ws1.Range("B5").Value = 5
$ws1.Range("C5").Value = 8
$ws1.Range("D5").Value = 13
$ws1.Range("E5").Value = 61.5

$ws1.Range("B6").Value = 2
$ws1.Range("C6").Value = 5
$ws1.Range("D6").Value = 7
$ws1.Range("E6").Value = 71.4

# --- Sheet 2: Total Accuracy ---
$ws2 = $wb.Worksheets.Item("Total Accuracy")

$ws2.Range("B2").Value = 7
$ws2.Range("C2").Value = 71
$ws2.Range("D2").Value = 78
$ws2.Range("E2").Value = 91

$ws2.Range("B3").Value = 2
$ws2.Range("C3").Value = 71
$ws2.Range("D3").Value = 73
$ws2.Range("E3").Value = 97.3

$ws2.Range("B4").Value = 5
$ws2.Range("C4").Value = 16
$ws2.Range("D4").Value = 21
$ws2.Range("E4").Value = 76.2

$ws2.Range("B5").Value = 5
$ws2.Range("C5").Value = 13
$ws2.Range("D5").Value = 18
$ws2.Range("E5").Value = 72.2

$ws2.Range("B6").Value = 1
$ws2.Range("C6").Value = 3
$ws2.Range("D6").Value = 4
$ws2.Range("E6").Value = 75

$wb.Save()
